# Add four new string metrics (MIN_STRING_CASE, MAX_STRING_CASE, MIN_STRING_ICASE,
# MAX_STRING_ICASE) as new rows 19-22 in the Profiler metric-mapping sheet,
# pushing the existing MAX/MIN/SUM/... rows down by four rows (now rows 23-34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four blank rows right before the current row 19 (where MAX used to live).
$ws.Range("A19:M22").Insert()

# The newly inserted rows don't inherit the table's direct formatting (borders/
# fill/alignment) from row insertion alone, so copy it over from the row above
# (row 18, which keeps the same per-column formatting used throughout the table).
$ws.Range("A18:M18").Copy()
$ws.Range("A19:M22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 19: MIN_STRING_CASE - applies to STRING and VARCHAR
$ws.Range("A19").Value = "MIN_STRING_CASE"
$ws.Range("I19").Value = "X"
$ws.Range("J19").Value = "X"

# Row 20: MAX_STRING_CASE - applies to STRING and VARCHAR
$ws.Range("A20").Value = "MAX_STRING_CASE"
$ws.Range("I20").Value = "X"
$ws.Range("J20").Value = "X"

# Row 21: MIN_STRING_ICASE - applies to STRING and VARCHAR
$ws.Range("A21").Value = "MIN_STRING_ICASE"
$ws.Range("I21").Value = "X"
$ws.Range("J21").Value = "X"

# Row 22: MAX_STRING_ICASE - applies to STRING and VARCHAR
$ws.Range("A22").Value = "MAX_STRING_ICASE"
$ws.Range("I22").Value = "X"
$ws.Range("J22").Value = "X"
